$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19
$ws.Cells.Item(19, 1).Value = "2024-10-05 14:06:44"
$ws.Cells.Item(19, 3).Value = 7
$ws.Cells.Item(19, 4).Value = 1
$ws.Cells.Item(19, 5).Value = 2
$ws.Cells.Item(19, 6).Value = 1
$ws.Cells.Item(19, 7).Value = 3
$ws.Cells.Item(19, 8).Value = 0
$ws.Cells.Item(19, 14).Value = 10
$ws.Cells.Item(19, 15).Value = 10
$ws.Cells.Item(19, 16).Value = 3
$ws.Cells.Item(19, 18).Value = 5
$ws.Cells.Item(19, 20).Value = 50
$ws.Cells.Item(19, 21).Value = 1
$ws.Cells.Item(19, 22).Value = "C:\Users\jonat\OneDrive\Escritorio\Repositorio\jonatha1992\Predictor_App\Data\Crupier.xlsx"
$ws.Cells.Item(19, 24).Value = "No es Simulación"
$ws.Cells.Item(19, 25).Value = 7

# Row 20
$ws.Cells.Item(20, 1).Value = "2024-10-05 14:54:13"
$ws.Cells.Item(20, 3).Value = 0
$ws.Cells.Item(20, 4).Value = 0
$ws.Cells.Item(20, 5).Value = 0
$ws.Cells.Item(20, 6).Value = 0
$ws.Cells.Item(20, 7).Value = 0
$ws.Cells.Item(20, 8).Value = 0
$ws.Cells.Item(20, 14).Value = 10
$ws.Cells.Item(20, 15).Value = 10
$ws.Cells.Item(20, 16).Value = 3
$ws.Cells.Item(20, 18).Value = 5
$ws.Cells.Item(20, 20).Value = 20
$ws.Cells.Item(20, 21).Value = 0
$ws.Cells.Item(20, 22).Value = "C:\Users\jonat\OneDrive\Escritorio\Repositorio\jonatha1992\Predictor_App\Data\Electromecanica.xlsx"
$ws.Cells.Item(20, 24).Value = "No es Simulación"
$ws.Cells.Item(20, 25).Value = 0
